# Generate Report for Handoff
# Adds a new row (row 3) for file "f63575dd-6795-4ae8-9a39-1119039d1969.md" to the
# "Overview", "zh-cn" and "de-de" worksheets / tables, mirroring the shape of the
# existing "59dea38f-1ea5-4f5b-b7be-8a5942b68881.md" row (row 2) on each sheet.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7871ddecdff83424bb9efbb88fa960fe74c57d08/e2e/"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> table "Overview" (columns A:G)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "f63575dd-6795-4ae8-9a39-1119039d1969.md"
$wsOverview.Range("B3").Value = "e2e\f63575dd-6795-4ae8-9a39-1119039d1969.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 10:47:45"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    ($repoBase + "f63575dd-6795-4ae8-9a39-1119039d1969.md"),
    "",
    "",
    "e2e\f63575dd-6795-4ae8-9a39-1119039d1969.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> table "zh-cn" (columns A:P)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "f63575dd-6795-4ae8-9a39-1119039d1969.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "f63575dd-6795-4ae8-9a39-1119039d1969.3dd23043a2bd596a4690c25cb133dc78e2c8b5a6.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 10:47:41"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    ($repoBase + "f63575dd-6795-4ae8-9a39-1119039d1969.md"),
    "",
    "",
    "f63575dd-6795-4ae8-9a39-1119039d1969.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" -> table "de-de" (columns A:P)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "f63575dd-6795-4ae8-9a39-1119039d1969.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "f63575dd-6795-4ae8-9a39-1119039d1969.3dd23043a2bd596a4690c25cb133dc78e2c8b5a6.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 10:47:45"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    ($repoBase + "f63575dd-6795-4ae8-9a39-1119039d1969.md"),
    "",
    "",
    "f63575dd-6795-4ae8-9a39-1119039d1969.md"
) | Out-Null
